$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("M3").Value = 1.11
$ws.Range("N3").Value = 6.5

# Row 7
$ws.Range("G7").Value = 2.15
$ws.Range("H7").Value = 2.92
$ws.Range("I7").Value = 3.5
$ws.Range("J7").Value = 2.75
$ws.Range("L7").Value = 4
$ws.Range("N7").Value = 6.8
$ws.Range("O7").Value = 1.4
$ws.Range("U7").Value = 1.85
$ws.Range("V7").Value = 1.75
$ws.Range("W7").Value = 6.4
$ws.Range("X7").Value = 9.75
$ws.Range("Z7").Value = 21
$ws.Range("AG7").Value = 8.75
$ws.Range("AH7").Value = 18
$ws.Range("AJ7").Value = 55
$ws.Range("AM7").Value = 800
$ws.Range("AN7").Value = 3.95
$ws.Range("AP7").Value = 19.5
$ws.Range("AR7").Value = 80
$ws.Range("AU7").Value = 6.9
$ws.Range("AV7").Value = 65
$ws.Range("AW7").Value = 5.3
$ws.Range("AY7").Value = 26

# Row 8
$ws.Range("G8").Value = 2.4
$ws.Range("M8").Value = 1.08
$ws.Range("N8").Value = 8
$ws.Range("S8").Value = 1.44
$ws.Range("T8").Value = 2.63
$ws.Range("AT8").Value = 2.63

# Row 13
$ws.Range("M13").Value = 1.02
$ws.Range("N13").Value = 7.1
